$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update names in column A
$ws.Range("A2").Value = "Nikki Klocko"
$ws.Range("A3").Value = "Sheldon McGlynn"
$ws.Range("A4").Value = "Yvette Mohr"
$ws.Range("A5").Value = "Candace Renner"

# Row 2: Applied/Accepted flags flip from 0 to 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

# Row 4: Applied/Accepted flags flip from 1 to 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
